$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text edits (reworded cells) ---------------------------------------
# A16: topic blurb gets expanded with links to the exercises/conditions
$ws.Range("A16").Value = "Задачи от темата: Преработка и постепенни промени, може да намерите задачите на: https://github.com/plamenna-petrova/Practical-Software-Development-11D/tree/master/src/Exercises/Refactoring-And-Progressive-Changes и условията: https://github.com/plamenna-petrova/Practical-Software-Development-11D/tree/master/src/Files/4.%20%D0%9F%D1%80%D0%B5%D1%80%D0%B0%D0%B1%D0%BE%D1%82%D0%BA%D0%B0%20%D0%B8%20%D0%BF%D0%BE%D1%81%D1%82%D0%B5%D0%BF%D0%B5%D0%BD%D0%BD%D0%B8%20%D0%BF%D1%80%D0%BE%D0%BC%D0%B5%D0%BD%D0%B8"

# A26: documentation-task reminder now numbered as task 4
$ws.Range("A26").Value = "4 задача. Да се документират всички задачи (техните методи и класове) чрез вградените средства на Visual Studio за XML документация"

# --- New column E: a 5th task-number column for rows 2-14 --------------
$ws.Range("E2").Value = 3
$ws.Range("E3").Value = 5
$ws.Range("E4").Value = 6
$ws.Range("E5").Value = 4
$ws.Range("E6").Value = 1
$ws.Range("E7").Value = 3
$ws.Range("E8").Value = 6
$ws.Range("E9").Value = 2
$ws.Range("E10").Value = 5
$ws.Range("E11").Value = 4
$ws.Range("E12").Value = 1
$ws.Range("E13").Value = 4
$ws.Range("E14").Value = 2

# Give the new column a left/right thin border (new style picked up by Excel)
$ws.Range("E2:E14").Borders.Item(7).LineStyle = 1   # xlEdgeLeft = 7, xlContinuous = 1
$ws.Range("E2:E14").Borders.Item(10).LineStyle = 1  # xlEdgeRight = 10

# --- Extend the conditional formatting color scales to include column E ---
$ws.Range("A1:D14,E2:E14").FormatConditions.Delete()
$cf1 = $ws.Range("A1:D14,E2:E14").FormatConditions.AddColorScale(3)
$cf1.ColorScaleCriteria.Item(1).FormatColor.Color = 7039595
$cf1.ColorScaleCriteria.Item(2).Type = 4
$cf1.ColorScaleCriteria.Item(2).FormatColor.Color = 16775676
$cf1.ColorScaleCriteria.Item(3).FormatColor.Color = 6530150

$cf2 = $ws.Range("A2:E14").FormatConditions.AddColorScale(3)
$cf2.ColorScaleCriteria.Item(1).FormatColor.Color = 7039595
$cf2.ColorScaleCriteria.Item(2).Type = 4
$cf2.ColorScaleCriteria.Item(2).FormatColor.Color = 16775676
$cf2.ColorScaleCriteria.Item(3).FormatColor.Color = 13012570

$cf3 = $ws.Range("E1").FormatConditions.AddColorScale(3)
$cf3.ColorScaleCriteria.Item(1).FormatColor.Color = 7039595
$cf3.ColorScaleCriteria.Item(2).Type = 4
$cf3.ColorScaleCriteria.Item(2).FormatColor.Color = 16775676
$cf3.ColorScaleCriteria.Item(3).FormatColor.Color = 6530150

# --- View state: scroll down and move the active selection -------------
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("B22").Select()
